# Update species list: always keep family infos, keep orig names also if not
# changed by GBIF.
#
# For several genus/family groups in column A, the bare genus/family name row
# (e.g. "Acer", "Apiaceae") is replaced by shifting the following species
# names up by one row, inserting a new "<Genus> species" row (Woodiness =
# "not found") at the alphabetically correct position among that genus's
# epithets, and leaving any remaining rows (already alphabetically after
# "species") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Acer (rows 3-7)
$ws.Range("A3").Value = "Acer campestre"
$ws.Range("B3").Value = "woody"
$ws.Range("A4").Value = "Acer negundo"
$ws.Range("A5").Value = "Acer platanoides"
$ws.Range("A6").Value = "Acer pseudoplatanus"
$ws.Range("A7").Value = "Acer species"
$ws.Range("B7").Value = "not found"

# Allium (row 11)
$ws.Range("A11").Value = "Allium species"
$ws.Range("B11").Value = "not found"

# Apiaceae (row 17)
$ws.Range("A17").Value = "Apiaceae species"

# Asteraceae (row 22)
$ws.Range("A22").Value = "Asteraceae species"
$ws.Range("B22").Value = "not found"

# Betula (row 27)
$ws.Range("A27").Value = "Betula species"

# Brassicaceae (row 32)
$ws.Range("A32").Value = "Brassicaceae species"
$ws.Range("B32").Value = "not found"

# Draba (row 66)
$ws.Range("A66").Value = "Draba species"
$ws.Range("B66").Value = "not found"

# Festuca (rows 77-80)
$ws.Range("A77").Value = "Festuca ovina"
$ws.Range("A78").Value = "Festuca rubra"
$ws.Range("A79").Value = "Festuca rupicola"
$ws.Range("A80").Value = "Festuca species"
$ws.Range("B80").Value = "not found"

# Geranium (rows 91-96)
$ws.Range("A91").Value = "Geranium molle"
$ws.Range("B91").Value = "herbaceous"
$ws.Range("A92").Value = "Geranium pratense"
$ws.Range("A93").Value = "Geranium pusillum"
$ws.Range("A94").Value = "Geranium pyrenaicum"
$ws.Range("B94").Value = "not found"
$ws.Range("A95").Value = "Geranium rotundifolium"
$ws.Range("B95").Value = "herbaceous"
$ws.Range("A96").Value = "Geranium species"
$ws.Range("B96").Value = "not found"

# Medicago (rows 123-125)
$ws.Range("A123").Value = "Medicago falcata"
$ws.Range("B123").Value = "herbaceous"
$ws.Range("A124").Value = "Medicago lupulina"
$ws.Range("A125").Value = "Medicago species"
$ws.Range("B125").Value = "not found"

# Poaceae (row 145)
$ws.Range("A145").Value = "Poaceae species"
$ws.Range("B145").Value = "not found"

# Prunus (rows 151-153)
$ws.Range("A151").Value = "Prunus avium"
$ws.Range("A152").Value = "Prunus mahaleb"
$ws.Range("A153").Value = "Prunus species"
$ws.Range("B153").Value = "not found"

# Rubus (rows 159-161)
$ws.Range("A159").Value = "Rubus caesius"
$ws.Range("B159").Value = "woody"
$ws.Range("A160").Value = "Rubus idaeus"
$ws.Range("A161").Value = "Rubus species"
$ws.Range("B161").Value = "not found"

# Senecio (rows 167-168)
$ws.Range("A167").Value = "Senecio jacobaea"
$ws.Range("B167").Value = "not found"
$ws.Range("A168").Value = "Senecio species"

# Triticum (row 202)
$ws.Range("A202").Value = "Triticum species"
$ws.Range("B202").Value = "not found"

# Vicia (rows 213-217)
$ws.Range("A213").Value = "Vicia cracca"
$ws.Range("B213").Value = "herbaceous"
$ws.Range("A214").Value = "Vicia hirsuta"
$ws.Range("A215").Value = "Vicia sativa"
$ws.Range("A216").Value = "Vicia sepium"
$ws.Range("B216").Value = "not assigned"
$ws.Range("A217").Value = "Vicia species"
$ws.Range("B217").Value = "not found"
